# Weekly driver report update for 2025-04-21
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bad Drivers section
$ws.Range("C3").Value = 432
$ws.Range("D3").Value = 82.2

$ws.Range("C4").Value = 1074
$ws.Range("D4").Value = 92.09999999999999

$ws.Range("C5").Value = 7260
$ws.Range("D5").Value = 92.59999999999999

$ws.Range("D6").Value = 97.3

$ws.Range("C7").Value = 77

# Totals row
$ws.Range("C8").Value = 8912

# Good Drivers section
$ws.Range("B16").Value = 449371
$ws.Range("B20").Value = 77999
